$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $row, $col, $text) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Extend formatted style (bold/border/centered) used by column A down to the
# new rows (11-13) before filling them in, by copying the format already
# used on an existing "A" cell.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A11:A13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 2314
$ws.Cells.Item(2, 3).Value = 95

# Row 3
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 2436
$ws.Cells.Item(3, 3).Value = 110

# Row 4
$ws.Cells.Item(4, 1).Value = 4
$ws.Cells.Item(4, 2).Value = 3035
$ws.Cells.Item(4, 3).Value = 193.5

# Row 5
$ws.Cells.Item(5, 1).Value = 5
$ws.Cells.Item(5, 2).Value = 3122
Set-TextValue $ws 5 3 "67.40"

# Row 6
$ws.Cells.Item(6, 1).Value = 6
$ws.Cells.Item(6, 2).Value = 3141
Set-TextValue $ws 6 3 "242.50"

# Row 7
$ws.Cells.Item(7, 1).Value = 8
$ws.Cells.Item(7, 2).Value = 3221
Set-TextValue $ws 7 3 "43.10"

# Row 8
$ws.Cells.Item(8, 1).Value = 10
$ws.Cells.Item(8, 2).Value = 3588
$ws.Cells.Item(8, 3).Value = 167

# Row 9
$ws.Cells.Item(9, 1).Value = 12
$ws.Cells.Item(9, 2).Value = 6104
Set-TextValue $ws 9 3 "179.00"

# Row 10
$ws.Cells.Item(10, 1).Value = 13
$ws.Cells.Item(10, 2).Value = 6138
Set-TextValue $ws 10 3 "214.50"

# Row 11 (new)
$ws.Cells.Item(11, 1).Value = 15
$ws.Cells.Item(11, 2).Value = 6271
$ws.Cells.Item(11, 3).Value = 300.5

# Row 12 (new)
$ws.Cells.Item(12, 1).Value = 16
$ws.Cells.Item(12, 2).Value = 6411
Set-TextValue $ws 12 3 "268.00"

# Row 13 (new)
$ws.Cells.Item(13, 1).Value = 19
$ws.Cells.Item(13, 2).Value = 8289
Set-TextValue $ws 13 3 "36.95"
